# Apply the price/volume updates from the latest cryptos.xlsx refresh.
# Column D ("Price") values are stored as plain text in the source data
# (e.g. "61.927.83"), so force Text format before writing each one to
# stop Excel from reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.868.52'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.410.69'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.26'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.49'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.636'
$ws.Range("E7").Value = '  +7.13%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.731'
$ws.Range("E9").Value = '  +5.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  +7.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.68'
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.06'
$ws.Range("E12").Value = '  +6.63%  '
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.950.87'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.21'
$ws.Range("E15").Value = '  +6.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000201'
$ws.Range("E16").Value = '  +38.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.404.14'
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.07'
$ws.Range("E18").Value = '  +5.50%  '
$ws.Range("E19").Value = '  +4.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '61.801.44'
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '453.62'
$ws.Range("E21").Value = '  +44.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '92.52'
$ws.Range("E22").Value = '  +10.31%  '
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.89'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.50'
$ws.Range("E26").Value = '  +12.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.75'
$ws.Range("E27").Value = '  +7.25%  '
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.56'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.74'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.98'
$ws.Range("E31").Value = '  +4.93%  '
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.80'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.167'
$ws.Range("E34").Value = '  -3.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0496'
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.30'
$ws.Range("E37").Value = '  +3.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("E39").Value = '  -1.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.134'
$ws.Range("E40").Value = '  +6.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.314'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.11'
$ws.Range("E43").Value = '  +1.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.18'
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("E46").Value = '  +8.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.51'
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.21'
$ws.Range("E48").Value = '  +3.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.759.63'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.105.60'
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '106.60'
$ws.Range("E51").Value = '  +27.29%  '
